$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style used by the rest of row 1 header cells (copy formatting
# from the neighboring header cell, same as the other columns use)
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data rows 2-53: Wins=66, Losses=96, Ties=0
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 66
    $ws.Cells.Item($r, 31).Value = 96
    $ws.Cells.Item($r, 32).Value = 0
}
